# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) values per-row with regenerated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 0
    8  = 3
    9  = 1
    10 = 1
    11 = 2
    12 = 0
    13 = 1
    14 = 2
    16 = 1
    17 = 1
    18 = 2
    20 = 1
    21 = 2
    22 = 0
    23 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
